# Update countries & provincias Spain
# - Add Ghana's updated figures (sorted position: between Costa de Marfil and Uruguay)
# - Add Barbados's updated figures (sorted position: between Etiopia and Jamaica)
# - Refresh Estados Unidos totals
# - Refresh Aruba totals
# - Refresh "last updated" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4): refresh totals ---
$ws.Range("B4").Value = 553057
$ws.Range("C4").Value = 20178
$ws.Range("E4").Value = 499968
$ws.Range("F4").Value = 11761
$ws.Range("G4").Value = 1143
$ws.Range("H4").Value = 21720

# --- Ghana: insert new row before Uruguay (row 93), drop the stale Ghana row ---
$ws.Rows.Item(93).Insert()
$ws.Range("A93").Value = "Ghana"
$ws.Range("B93").Value = 566
$ws.Range("C93").Value = 158
$ws.Range("D93").Value = 4
$ws.Range("E93").Value = 554
$ws.Range("F93").Value = 2
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 8
# The old Ghana row has shifted down to row 98 (Costa de Marfil 92, Ghana(new) 93,
# Uruguay 94, Burkina Faso 95, Niger 96, Albania 97, Ghana(old, stale) 98) - remove it.
$ws.Rows.Item(98).Delete()

# --- Barbados: insert new row before Jamaica (row 140), drop the stale Barbados row ---
$ws.Rows.Item(140).Insert()
$ws.Range("A140").Value = "Barbados"
$ws.Range("B140").Value = 71
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 11
$ws.Range("E140").Value = 56
$ws.Range("F140").Value = 4
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 4
# The old Barbados row has shifted down to row 142 (Etiopia 139, Barbados(new) 140,
# Jamaica 141, Barbados(old, stale) 142, Congo 143) - remove it.
$ws.Rows.Item(142).Delete()

# --- Aruba: refresh active/recovered/critical counts ---
$ws.Range("D135").Value = 32
$ws.Range("E135").Value = 60
$ws.Range("F135").Value = 1

# --- Footer: refresh "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 21:52"
